$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '33.955.21'
Set-TextValue $ws.Range("E2") '  +9.12%  '
Set-TextValue $ws.Range("D3") '1.783.66'
Set-TextValue $ws.Range("E3") '  +5.91%  '
Set-TextValue $ws.Range("D4") '1.00'
Set-TextValue $ws.Range("E4") '  +0.34%  '
Set-TextValue $ws.Range("D5") '224.77'
Set-TextValue $ws.Range("E5") '  +2.13%  '
Set-TextValue $ws.Range("D6") '0.556'
Set-TextValue $ws.Range("E6") '  +4.15%  '
Set-TextValue $ws.Range("D7") '1.00'
Set-TextValue $ws.Range("E7") '  +0.45%  '
Set-TextValue $ws.Range("D8") '30.66'
Set-TextValue $ws.Range("E8") '  +3.97%  '
Set-TextValue $ws.Range("D9") '46.33'
Set-TextValue $ws.Range("E9") '  +3.04%  '
Set-TextValue $ws.Range("E10") '  +4.51%  '
Set-TextValue $ws.Range("D11") '0.0660'
Set-TextValue $ws.Range("E11") '  +3.35%  '
Set-TextValue $ws.Range("E12") '  +2.06%  '
Set-TextValue $ws.Range("D13") '2.041.88'
Set-TextValue $ws.Range("E13") '  +6.10%  '
Set-TextValue $ws.Range("D14") '1.786.62'
Set-TextValue $ws.Range("E14") '  +6.13%  '
Set-TextValue $ws.Range("D15") '0.626'
Set-TextValue $ws.Range("E15") '  +2.66%  '
Set-TextValue $ws.Range("D16") '33.932.48'
Set-TextValue $ws.Range("E16") '  +9.07%  '
Set-TextValue $ws.Range("D17") '9.97'
Set-TextValue $ws.Range("E17") '  -3.78%  '
Set-TextValue $ws.Range("D18") '4.19'
Set-TextValue $ws.Range("E18") '  +2.04%  '
Set-TextValue $ws.Range("D19") '68.47'
Set-TextValue $ws.Range("E19") '  +2.62%  '
Set-TextValue $ws.Range("D20") '251.50'
Set-TextValue $ws.Range("E20") '  +1.20%  '
Set-TextValue $ws.Range("E21") '  +2.67%  '
Set-TextValue $ws.Range("D22") '0.999'
Set-TextValue $ws.Range("E22") '  +0.25%  '
Set-TextValue $ws.Range("E23") '  +2.32%  '
Set-TextValue $ws.Range("D24") '4.21'
Set-TextValue $ws.Range("E24") '  -1.50%  '
Set-TextValue $ws.Range("E25") '  -0.30%  '
Set-TextValue $ws.Range("D26") '158.24'
Set-TextValue $ws.Range("E26") '  -0.23%  '
Set-TextValue $ws.Range("D27") '16.44'
Set-TextValue $ws.Range("E27") '  +3.13%  '
Set-TextValue $ws.Range("E28") '  +1.76%  '
Set-TextValue $ws.Range("E29") '  +3.45%  '
Set-TextValue $ws.Range("D30") '1.00'
Set-TextValue $ws.Range("E30") '  +0.33%  '
Set-TextValue $ws.Range("E31") '  +7.91%  '
Set-TextValue $ws.Range("D32") '0.0511'
Set-TextValue $ws.Range("E32") '  +2.94%  '
Set-TextValue $ws.Range("E33") '  +4.34%  '
Set-TextValue $ws.Range("E34") '  +5.57%  '
Set-TextValue $ws.Range("D35") '1.495.28'
Set-TextValue $ws.Range("E35") '  -1.48%  '
Set-TextValue $ws.Range("E36") '  +0.50%  '
Set-TextValue $ws.Range("E37") '  +3.05%  '
Set-TextValue $ws.Range("D38") '0.631'
Set-TextValue $ws.Range("E38") '  +3.34%  '
Set-TextValue $ws.Range("E39") '  +2.93%  '
Set-TextValue $ws.Range("D40") '82.98'
Set-TextValue $ws.Range("E40") '  -2.27%  '
Set-TextValue $ws.Range("D41") '2.36'
Set-TextValue $ws.Range("E41") '  +2.98%  '
Set-TextValue $ws.Range("D42") '2.70'
Set-TextValue $ws.Range("E42") '  +0.97%  '
Set-TextValue $ws.Range("D43") '0.885'
Set-TextValue $ws.Range("E43") '  +4.89%  '
Set-TextValue $ws.Range("E44") '  +2.08%  '
Set-TextValue $ws.Range("E45") '  +1.38%  '
Set-TextValue $ws.Range("D46") '1.06'
Set-TextValue $ws.Range("E46") '  +2.87%  '
Set-TextValue $ws.Range("D47") '1.938.66'
Set-TextValue $ws.Range("E47") '  +6.79%  '
Set-TextValue $ws.Range("D48") '5.75'
Set-TextValue $ws.Range("E48") '  +3.00%  '
Set-TextValue $ws.Range("E49") '  +0.33%  '
Set-TextValue $ws.Range("D50") '11.88'
Set-TextValue $ws.Range("E50") '  +12.87%  '
Set-TextValue $ws.Range("D51") '50.85'
Set-TextValue $ws.Range("E51") '  -3.16%  '
